$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a fresh data row at row 77 (pushing the existing
# rows 77-99 down to 78-100, carrying their original data and formatting
# along with them).
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with this week's new observation.
$ws.Cells.Item(77, 1).Value = 10
$ws.Cells.Item(77, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(77, 3).Value = "La Araucanía"
$ws.Cells.Item(77, 4).Value = 44754
$ws.Cells.Item(77, 5).Value = 9
$ws.Cells.Item(77, 6).Value = 100114002
$ws.Cells.Item(77, 7).Value = "Camote"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 30
$ws.Cells.Item(77, 11).Value = 20000
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = 20000
$ws.Cells.Item(77, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(77, 15).Value = "Perú"
$ws.Cells.Item(77, 16).Value = 1000
$ws.Cells.Item(77, 17).Value = 20
$ws.Cells.Item(77, 18).Value = "Hortaliza"
